$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vocabulary rows to append: word, definition, example 1, example 2
$newRows = @(
    @{ Row = 103; Height = 75;  Word = "severity";  Def = "severe problems, injuries, illnesses etc are very bad or very serious"; Ex1 = "he risk and severity of sunborn depend on he body's natural skin color."; Ex2 = "His injuries were quite severe." },
    @{ Row = 104; Height = 60;  Word = "sensation"; Def = "a feeling that you get from one of your five senses, especially the sense of touch"; Ex1 = "I experienced no sensation in my left foot."; Ex2 = "One sign of a heart attack is a tingling sensation in the left arm." },
    @{ Row = 105; Height = 45;  Word = "smuggle";   Def = "to take something or someone illegally from one country to another"; Ex1 = "if you try to smuggle drug you are stupid."; Ex2 = "The guns were smuggled across the border." },
    @{ Row = 106; Height = 75;  Word = "slope";     Def = "a surface of which one end or side is at a higher level than another; a rising or falling surface."; Ex1 = "the house builders slopped the roof..."; Ex2 = "the roof should have a slope sufficient for proper drainage" },
    @{ Row = 107; Height = 105; Word = "soak";      Def = "if you soak something, or if you let it soak, you keep it covered with a liquid for a period of time, especially in order to make it softer or easier to clean"; Ex1 = "Soak the clothes in cold water."; Ex2 = "soak the beans overnight in water" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Word
    $ws.Cells.Item($rowNum, 2).Value = $r.Def
    $ws.Cells.Item($rowNum, 3).Value = $r.Ex1
    $ws.Cells.Item($rowNum, 4).Value = $r.Ex2

    $rowRange = $ws.Range("A" + $rowNum + ":D" + $rowNum)
    $rowRange.WrapText = $true
    $rowRange.VerticalAlignment = -4160  # xlTop, matches style used by the rest of the table body

    $ws.Rows.Item($rowNum).RowHeight = $r.Height
}

# Reflect where the user ended up after typing the last entry: the view
# had scrolled down so row 101 is at the top, with the cursor past the
# last filled cell (E107).
try {
    $excel.ActiveWindow.ScrollRow = 101
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scrolling the viewport isn't always settable in every host; ignore.
}
$ws.Range("E107").Select() | Out-Null
